# Apply the edits described in the commit: update elapsed/age boundary
# values in columns C/D for a handful of rows, turning several of the
# column-C entries into "= previous row's D" formulas (as Excel would
# produce when you fill a formula down over existing rows), and fix one
# mislabeled phoneme in F336.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Block 1: rows 278-284
# ---------------------------------------------------------------------
$ws.Range("C278").Value = 127.1
$ws.Range("D278").Value = 127.25

$ws.Range("C279").Formula = "=D278"
$ws.Range("D279").Value = 127.35

# C280:C282 fill down as one shared-formula group (=D(row-1))
$ws.Range("C280:C282").Formula = "=D279"
$ws.Range("D280").Value = 127.5
$ws.Range("D281").Value = 127.65
$ws.Range("D282").Value = 127.85

$ws.Range("C283").Formula = "=D282"
$ws.Range("D283").Value = 128.19999999999999

$ws.Range("C284").Formula = "=D283"
$ws.Range("D284").Value = 128.63

# ---------------------------------------------------------------------
# Block 2: rows 297-307
# ---------------------------------------------------------------------
$ws.Range("C297").Value = 133.34
$ws.Range("D297").Value = 133.5

$ws.Range("C298").Value = 133.5
$ws.Range("D298").Value = 133.69999999999999

$ws.Range("D299").Value = 124.6

$ws.Range("C300").Value = 124.6
$ws.Range("D300").Value = 125

$ws.Range("C301").Formula = "=D300"

# C302:C307 fill down as one shared-formula group (=D(row-1))
$ws.Range("C302:C307").Formula = "=D301"

# ---------------------------------------------------------------------
# Block 3: rows 325-336
# ---------------------------------------------------------------------
$ws.Range("C325").Value = 145.035
$ws.Range("D325").Value = 145.63999999999999

$ws.Range("C326").Formula = "=D325"
$ws.Range("D326").Value = 145.72

$ws.Range("C327").Formula = "=D326"
$ws.Range("D327").Value = 146.46

$ws.Range("C328").Formula = "=D327"
$ws.Range("D328").Value = 147.1

# C329:C336 fill down as one shared-formula group (=D(row-1))
$ws.Range("C329:C336").Formula = "=D328"
$ws.Range("D329").Value = 147.07
$ws.Range("D330").Value = 147.32
$ws.Range("D331").Value = 147.33500000000001
$ws.Range("D332").Value = 148.465

# F336: relabel phoneme "i" -> "n"
$ws.Range("F336").Value = "n"

# ---------------------------------------------------------------------
# Sheet view / selection housekeeping (matches the author scrolling to
# the bottom of the sheet and selecting F337 while zoomed to 139%).
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 139
$ws.Range("F337").Select()
